$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "DMG1012T-7"
$ws.Range("C6").Value = "DMG1012T-7DICT-ND"

$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null

$ws.Range("Y3").Select() | Out-Null
